$wb = $excel.ActiveWorkbook

# --- Sheet 1: birth details ---
$ws1 = $wb.Worksheets.Item("Sheet 1")

$ws1.Range("B4").Value = "19:34:00"
$ws1.Range("B5").Value = "Bantwal,Karnataka"
$ws1.Range("B6").Value = 12.8953789
$ws1.Range("B7").Value = 75.0408618

# --- Sheet 2: planetary positions ---
$ws2 = $wb.Worksheets.Item("Sheet 2")

# Row 3 - Ascendant
$ws2.Range("B3").Value = "Leo"
$ws2.Range("C3").Value = "Sun"
$ws2.Range("D3").Value = "Magha"
$ws2.Range("E3").Value = "Ketu"
$ws2.Range("F3").Value = 2.818248223016099

# Row 4 - Sun
$ws2.Range("F4").Value = 286.9412862534865
$ws2.Range("J4").Value = 6

# Row 5 - Moon
$ws2.Range("F5").Value = 12.57328065113894
$ws2.Range("J5").Value = 9

# Row 6 - Mercury
$ws2.Range("F6").Value = 291.8539258992433
$ws2.Range("J6").Value = 6

# Row 7 - Venus
$ws2.Range("D7").Value = "Uttara Bhadrapada"
$ws2.Range("E7").Value = "Saturn"
$ws2.Range("F7").Value = 333.5552974085365
$ws2.Range("J7").Value = 8

# Row 8 - Mars
$ws2.Range("F8").Value = 76.51867948886958
$ws2.Range("J8").Value = 11

# Row 9 - Jupiter
$ws2.Range("F9").Value = 170.9261886981888
$ws2.Range("J9").Value = 2

# Row 10 - Saturn
$ws2.Range("F10").Value = 295.952888004387
$ws2.Range("J10").Value = 6

# Row 11 - Uranus
$ws2.Range("F11").Value = 265.6412131325364
$ws2.Range("J11").Value = 5

# Row 12 - Neptune
$ws2.Range("F12").Value = 265.6974764303666
$ws2.Range("J12").Value = 5

# Row 13 - Pluto
$ws2.Range("F13").Value = 211.5399106679366
$ws2.Range("J13").Value = 4

# Row 14 - Rahu
$ws2.Range("F14").Value = 235.0936835767791
$ws2.Range("J14").Value = 4

# Row 15 - Ketu
$ws2.Range("F15").Value = 55.09368357677914
$ws2.Range("J15").Value = 10
